# Adding breathcounter and shared_table example forms
#
# Concretely (per the target diff) this inserts a new "geopoint" survey
# question -- a demo of the ODK "Capture your location" intent -- right
# after the existing barcode-scanning demo row and before the
# picture-taking demo row on the "survey" sheet. Inserting the row shifts
# every subsequent row down by one; Excel re-points the shared-string
# table and row styles automatically, so no other sheet needs touching.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Push rows 9-39 down to 10-40, opening up a blank row 9.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row with the geopoint example question.
$ws.Cells.Item(9, 1).Value = "geopoint"
$ws.Cells.Item(9, 3).Value = "geopoint"
$ws.Cells.Item(9, 4).Value = "Capture your location"
